$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).NumberFormat = "@"
$ws.Cells.Item(2,4).Value = "63.591.52"
$ws.Cells.Item(2,5).Value = "  +0.23%  "
$ws.Cells.Item(3,4).NumberFormat = "@"
$ws.Cells.Item(3,4).Value = "3.349.68"
$ws.Cells.Item(3,5).Value = "  +2.09%  "
$ws.Cells.Item(4,5).Value = "  +0.02%  "
$ws.Cells.Item(5,4).NumberFormat = "@"
$ws.Cells.Item(5,4).Value = "531.75"
$ws.Cells.Item(5,5).Value = "  +3.02%  "
$ws.Cells.Item(6,4).NumberFormat = "@"
$ws.Cells.Item(6,4).Value = "172.77"
$ws.Cells.Item(6,5).Value = "  -4.82%  "
$ws.Cells.Item(7,5).Value = "  -0.10%  "
$ws.Cells.Item(8,4).NumberFormat = "@"
$ws.Cells.Item(8,4).Value = "3.348.72"
$ws.Cells.Item(8,5).Value = "  +1.82%  "
$ws.Cells.Item(9,5).Value = "  +0.03%  "
$ws.Cells.Item(10,5).Value = "  -1.67%  "
$ws.Cells.Item(11,4).NumberFormat = "@"
$ws.Cells.Item(11,4).Value = "53.13"
$ws.Cells.Item(11,5).Value = "  -9.82%  "
$ws.Cells.Item(12,5).Value = "  +1.74%  "
$ws.Cells.Item(13,5).Value = "  +0.54%  "
$ws.Cells.Item(14,4).NumberFormat = "@"
$ws.Cells.Item(14,4).Value = "9.18"
$ws.Cells.Item(14,5).Value = "  +0.70%  "
$ws.Cells.Item(15,4).NumberFormat = "@"
$ws.Cells.Item(15,4).Value = "3.876.75"
$ws.Cells.Item(15,5).Value = "  +1.93%  "
$ws.Cells.Item(16,4).NumberFormat = "@"
$ws.Cells.Item(16,4).Value = "3.345.28"
$ws.Cells.Item(16,5).Value = "  +1.98%  "
$ws.Cells.Item(17,4).NumberFormat = "@"
$ws.Cells.Item(17,4).Value = "0.117"
$ws.Cells.Item(17,5).Value = "  -0.71%  "
$ws.Cells.Item(18,4).NumberFormat = "@"
$ws.Cells.Item(18,4).Value = "17.49"
$ws.Cells.Item(18,5).Value = "  -0.61%  "
$ws.Cells.Item(19,4).NumberFormat = "@"
$ws.Cells.Item(19,4).Value = "63.548.49"
$ws.Cells.Item(19,5).Value = "  +0.38%  "
$ws.Cells.Item(20,5).Value = "  +2.36%  "
$ws.Cells.Item(21,4).NumberFormat = "@"
$ws.Cells.Item(21,4).Value = "0.966"
$ws.Cells.Item(21,5).Value = "  +2.10%  "
$ws.Cells.Item(22,4).NumberFormat = "@"
$ws.Cells.Item(22,4).Value = "372.70"
$ws.Cells.Item(22,5).Value = "  +0.40%  "
$ws.Cells.Item(23,4).NumberFormat = "@"
$ws.Cells.Item(23,4).Value = "4.14"
$ws.Cells.Item(23,5).Value = "  +6.56%  "
$ws.Cells.Item(24,2).Value = "RenderToken"
$ws.Cells.Item(24,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(24,4).NumberFormat = "@"
$ws.Cells.Item(24,4).Value = "11.29"
$ws.Cells.Item(24,5).Value = "  +0.41%  "
$ws.Cells.Item(25,2).Value = "Litecoin"
$ws.Cells.Item(25,3).Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Cells.Item(25,4).NumberFormat = "@"
$ws.Cells.Item(25,4).Value = "81.64"
$ws.Cells.Item(25,5).Value = "  +1.83%  "
$ws.Cells.Item(26,5).Value = "  +2.66%  "
$ws.Cells.Item(27,4).NumberFormat = "@"
$ws.Cells.Item(27,4).Value = "6.18"
$ws.Cells.Item(27,5).Value = "  +3.68%  "
$ws.Cells.Item(28,4).NumberFormat = "@"
$ws.Cells.Item(28,4).Value = "2.69"
$ws.Cells.Item(28,5).Value = "  +1.86%  "
$ws.Cells.Item(29,4).NumberFormat = "@"
$ws.Cells.Item(29,4).Value = "11.33"
$ws.Cells.Item(30,4).NumberFormat = "@"
$ws.Cells.Item(30,4).Value = "8.28"
$ws.Cells.Item(30,5).Value = "  -0.38%  "
$ws.Cells.Item(31,4).NumberFormat = "@"
$ws.Cells.Item(31,4).Value = "28.87"
$ws.Cells.Item(31,5).Value = "  +1.29%  "
$ws.Cells.Item(32,4).NumberFormat = "@"
$ws.Cells.Item(32,4).Value = "645.52"
$ws.Cells.Item(32,5).Value = "  -0.42%  "
$ws.Cells.Item(33,4).NumberFormat = "@"
$ws.Cells.Item(33,4).Value = "6.42"
$ws.Cells.Item(33,5).Value = "  -3.97%  "
$ws.Cells.Item(34,4).NumberFormat = "@"
$ws.Cells.Item(34,4).Value = "11.20"
$ws.Cells.Item(34,5).Value = "  +0.36%  "
$ws.Cells.Item(35,5).Value = "  +1.60%  "
$ws.Cells.Item(36,4).NumberFormat = "@"
$ws.Cells.Item(36,4).Value = "58.05"
$ws.Cells.Item(36,5).Value = "  -2.13%  "
$ws.Cells.Item(37,5).Value = "  +0.03%  "
$ws.Cells.Item(38,4).NumberFormat = "@"
$ws.Cells.Item(38,4).Value = "36.97"
$ws.Cells.Item(38,5).Value = "  +2.63%  "
$ws.Cells.Item(39,4).NumberFormat = "@"
$ws.Cells.Item(39,4).Value = "0.381"
$ws.Cells.Item(39,5).Value = "  -1.52%  "
$ws.Cells.Item(40,4).NumberFormat = "@"
$ws.Cells.Item(40,4).Value = "0.0₃0728"
$ws.Cells.Item(40,5).Value = "  +12.23%  "
$ws.Cells.Item(41,5).Value = "  +0.11%  "
$ws.Cells.Item(42,2).Value = "Fetch.AI"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "2.61"
$ws.Cells.Item(42,5).Value = "  +7.97%  "
$ws.Cells.Item(43,2).Value = "Kaspa"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "0.125"
$ws.Cells.Item(43,5).Value = "  -0.09%  "
$ws.Cells.Item(44,4).NumberFormat = "@"
$ws.Cells.Item(44,4).Value = "2.929.28"
$ws.Cells.Item(44,5).Value = "  -1.20%  "
$ws.Cells.Item(45,4).NumberFormat = "@"
$ws.Cells.Item(45,4).Value = "2.99"
$ws.Cells.Item(45,5).Value = "  +7.12%  "
$ws.Cells.Item(46,5).Value = "  +3.14%  "
$ws.Cells.Item(47,4).NumberFormat = "@"
$ws.Cells.Item(47,4).Value = "0.0397"
$ws.Cells.Item(47,5).Value = "  +2.77%  "
$ws.Cells.Item(48,5).Value = "  -2.37%  "
$ws.Cells.Item(49,5).Value = "  +3.69%  "
$ws.Cells.Item(50,2).Value = "Stellar"
$ws.Cells.Item(50,3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(50,4).NumberFormat = "@"
$ws.Cells.Item(50,4).Value = "0.125"
$ws.Cells.Item(50,5).Value = "  -0.25%  "
$ws.Cells.Item(51,2).Value = "Monero"
$ws.Cells.Item(51,3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(51,4).NumberFormat = "@"
$ws.Cells.Item(51,4).Value = "136.71"
$ws.Cells.Item(51,5).Value = "  +4.38%  "
